$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rushing")
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "J.Hurts"
$ws.Cells.Item(2,3).Value = 50
$ws.Cells.Item(2,4).Value = 44
$ws.Cells.Item(2,5).Value = 35
$ws.Cells.Item(2,6).Value = 30
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "G.Minshew"
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "M.Sanders"
$ws.Cells.Item(4,3).Value = 124
$ws.Cells.Item(4,4).Value = 84
$ws.Cells.Item(4,5).Value = 13
$ws.Cells.Item(4,6).Value = 30
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "B.Scott"
$ws.Cells.Item(5,3).Value = 51
$ws.Cells.Item(5,4).Value = 17
$ws.Cells.Item(5,5).Value = 13
$ws.Cells.Item(5,6).Value = 15
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "K.Gainwell"
$ws.Cells.Item(6,3).Value = 7
$ws.Cells.Item(6,4).Value = 5
$ws.Cells.Item(6,5).Value = 7
$ws.Cells.Item(6,6).Value = 5
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "J.Howard"
$ws.Cells.Item(7,3).Value = 2
$ws.Cells.Item(7,4).Value = 5
$ws.Cells.Item(7,5).Value = 10
$ws.Cells.Item(7,6).Value = 5
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "J.Reagor"
$ws.Cells.Item(8,3).Value = 5
$ws.Cells.Item(8,4).Value = 1
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "Q.Watkins"
$ws.Cells.Item(9,3).Value = 1
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "M.Walker"
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(9,1).Copy()
$ws.Cells.Item(10,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws = $wb.Worksheets.Item("Receiving")
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "M.Sanders"
$ws.Cells.Item(2,3).Value = 44
$ws.Cells.Item(2,4).Value = 32
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 8
$ws.Cells.Item(2,8).Value = 7
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "B.Scott"
$ws.Cells.Item(3,3).Value = 13
$ws.Cells.Item(3,4).Value = 10
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 1
$ws.Cells.Item(3,8).Value = 0
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "K.Gainwell"
$ws.Cells.Item(4,3).Value = 13
$ws.Cells.Item(4,4).Value = 11
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 1
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "J.Howard"
$ws.Cells.Item(5,3).Value = 1
$ws.Cells.Item(5,4).Value = 1
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "D.Smith"
$ws.Cells.Item(6,3).Value = 72
$ws.Cells.Item(6,4).Value = 50
$ws.Cells.Item(6,5).Value = 39
$ws.Cells.Item(6,6).Value = 16
$ws.Cells.Item(6,7).Value = 10
$ws.Cells.Item(6,8).Value = 5
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "J.Reagor"
$ws.Cells.Item(7,3).Value = 45
$ws.Cells.Item(7,4).Value = 28
$ws.Cells.Item(7,5).Value = 12
$ws.Cells.Item(7,6).Value = 4
$ws.Cells.Item(7,7).Value = 5
$ws.Cells.Item(7,8).Value = 3
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "Q.Watkins"
$ws.Cells.Item(8,3).Value = 37
$ws.Cells.Item(8,4).Value = 28
$ws.Cells.Item(8,5).Value = 24
$ws.Cells.Item(8,6).Value = 11
$ws.Cells.Item(8,7).Value = 9
$ws.Cells.Item(8,8).Value = 5
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "G.Ward"
$ws.Cells.Item(9,3).Value = 9
$ws.Cells.Item(9,4).Value = 5
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 7
$ws.Cells.Item(9,8).Value = 3
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "J.Arcega-Whiteside"
$ws.Cells.Item(10,3).Value = 1
$ws.Cells.Item(10,4).Value = 1
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = 0
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "D.Goedert"
$ws.Cells.Item(11,3).Value = 64
$ws.Cells.Item(11,4).Value = 44
$ws.Cells.Item(11,5).Value = 23
$ws.Cells.Item(11,6).Value = 17
$ws.Cells.Item(11,7).Value = 7
$ws.Cells.Item(11,8).Value = 5
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "J.Stoll"
$ws.Cells.Item(12,3).Value = 4
$ws.Cells.Item(12,4).Value = 3
$ws.Cells.Item(12,5).Value = 0
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 1
$ws.Cells.Item(12,8).Value = 1
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "T.Jackson"
$ws.Cells.Item(13,3).Value = 1
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 0
$ws.Cells.Item(13,8).Value = 0
$ws.Cells.Item(12,1).Copy()
$ws.Cells.Item(13,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
